# shopping cart page complete module push
# Populate column A (rows 2-39) of Sheet1 with the quick-cart example data
# and move the active selection to M10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    39015,
    7907,
    73189,
    102914,
    7645,
    82832,
    74999,
    5202,
    7260,
    59168,
    77545,
    2576,
    6173,
    88349,
    28014,
    78968,
    59555,
    29381,
    34590,
    79814,
    100236,
    71564,
    78458,
    53828,
    99667,
    106332,
    26366,
    27766,
    82720,
    80986,
    85328,
    35876,
    76508,
    98547,
    87112,
    7063,
    86000,
    9156
)

$row = 2
foreach ($v in $values) {
    $ws.Cells.Item($row, 1).Value = $v
    $row = $row + 1
}

# Move / record the active selection as in the target workbook
$null = $ws.Range("M10").Select()
